$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (shifts old C..O to D..P).
$ws.Range("C:C").Insert()

# The newly inserted column C currently carries formatting copied from the
# old column B (left neighbour). The source workbook instead keeps the
# formatting that used to belong to column C (now shifted to D), so copy
# the formats back from D into C to match.
$ws.Range("D1:D16").Copy()
$ws.Range("C1:C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header + unit row for the new "T_evap [deg C]" column.
$ws.Range("C1").Value = "T_evap"
$ws.Range("C2").Value = "[deg C]"

# New T_evap data values (row 3 .. row 16).
$ws.Range("C3").Value = 9.2899999999999991
$ws.Range("C4").Value = 4.66
$ws.Range("C5").Value = 10.36
$ws.Range("C6").Value = 0.03
$ws.Range("C7").Value = 4.9000000000000004
$ws.Range("C8").Value = 10.49
$ws.Range("C9").Value = -5.84
$ws.Range("C10").Value = 0.65
$ws.Range("C11").Value = 5.74
$ws.Range("C12").Value = -9.3000000000000007
$ws.Range("C13").Value = -5.351
$ws.Range("C14").Value = -9.89
$ws.Range("C15").Value = -15.01
$ws.Range("C16").Value = -22.71

# Move the selection/active cell to G2 (matches the edited file's saved view).
$ws.Range("G2").Select()
